$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (D) and "Volume(1h)" (E) columns with the latest crypto snapshot values.
# Each row keeps its text formatting (e.g. trailing zeros, thousands separators, padded
# percentages) by forcing a text number format before writing the value.
$updates = @(
    @{ Row = 2; D = "66.611.89"; E = "  +0.12%  " }
    @{ Row = 3; D = "3.521.37"; E = "  -2.18%  " }
    @{ Row = 4; D = "0.999"; E = "  -0.08%  " }
    @{ Row = 5; D = "607.54"; E = "  -0.29%  " }
    @{ Row = 6; D = "143.46"; E = "  -3.87%  " }
    @{ Row = 7; D = "3.520.90"; E = "  -2.18%  " }
    @{ Row = 8; D = "0.999"; E = "  -0.17%  " }
    @{ Row = 9; D = "0.510"; E = "  +4.08%  " }
    @{ Row = 10; D = "7.75"; E = "  -3.67%  " }
    @{ Row = 11; D = $null; E = "  -4.34%  " }
    @{ Row = 12; D = "0.409"; E = "  -1.71%  " }
    @{ Row = 13; D = "4.117.22"; E = "  -2.21%  " }
    @{ Row = 14; D = "0.0000195"; E = "  -6.47%  " }
    @{ Row = 15; D = "28.70"; E = "  -3.91%  " }
    @{ Row = 16; D = "3.524.19"; E = "  -2.05%  " }
    @{ Row = 17; D = $null; E = "  +0.40%  " }
    @{ Row = 18; D = "66.491.85"; E = "  -0.16%  " }
    @{ Row = 19; D = "10.78"; E = "  -6.89%  " }
    @{ Row = 20; D = "6.16"; E = "  -3.49%  " }
    @{ Row = 21; D = "14.69"; E = "  -2.86%  " }
    @{ Row = 22; D = "423.46"; E = "  -0.91%  " }
    @{ Row = 23; D = "0.589"; E = "  -4.96%  " }
    @{ Row = 24; D = "76.88"; E = "  -2.39%  " }
    @{ Row = 25; D = "3.660.34"; E = "  -2.24%  " }
    @{ Row = 26; D = $null; E = "  +0.11%  " }
    @{ Row = 27; D = "0.0000114"; E = "  -6.03%  " }
    @{ Row = 28; D = "7.93"; E = "  -4.76%  " }
    @{ Row = 29; D = "2.47"; E = "  -2.06%  " }
    @{ Row = 30; D = "8.94"; E = "  -5.25%  " }
    @{ Row = 31; D = "1.00"; E = $null }
    @{ Row = 32; D = "3.526.37"; E = "  -1.96%  " }
    @{ Row = 33; D = "0.155"; E = "  -1.66%  " }
    @{ Row = 34; D = "24.25"; E = "  -4.76%  " }
    @{ Row = 36; D = "1.34"; E = "  -9.04%  " }
    @{ Row = 37; D = "7.59"; E = "  -3.42%  " }
    @{ Row = 38; D = "1.63"; E = "  -4.00%  " }
    @{ Row = 39; D = "173.68"; E = "  -2.20%  " }
    @{ Row = 40; D = "5.23"; E = "  -7.53%  " }
    @{ Row = 41; D = $null; E = "  -4.94%  " }
    @{ Row = 42; D = "4.99"; E = "  -4.88%  " }
    @{ Row = 43; D = "0.855"; E = "  -4.85%  " }
    @{ Row = 44; D = "45.49"; E = "  -0.87%  " }
    @{ Row = 45; D = "1.77"; E = "  -7.00%  " }
    @{ Row = 46; D = "1.00"; E = "  +0.03%  " }
    @{ Row = 47; D = "2.36"; E = "  -8.32%  " }
    @{ Row = 48; D = "7.07"; E = "  -1.84%  " }
    @{ Row = 49; D = "1.12"; E = "  -5.23%  " }
    @{ Row = 50; D = "22.73"; E = "  -5.57%  " }
    @{ Row = 51; D = "0.907"; E = "  -4.89%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($u.E -ne $null) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.ClearFormats()
    }
}
